# Update the build/version string throughout the workbook.
# Old: mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)
# New: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..." line
$wsAbout.Range("A2").Value2 = "Version: " + $newVersion

# A6: Recommended Citation line, with the version embedded in single quotes
$wsAbout.Range("A6").Value2 = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Broadmeadow Coal Mine, Australia, M0016, version ' + "'" + $newVersion + "'" + ". (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Column S (build_version) rows 2 through 33 on the "Boundaries and methane sources" sheet
for ($r = 2; $r -le 33; $r++) {
    $cell = $wsBoundaries.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value2 = $newVersion
    }
}
